$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D34').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D35').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D114').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D122').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D157').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D161').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D164').Value = 'Downgrade Since Last Fixed Quarter Date'
$ws.Range('D192').Value = 'Downgrade Since Last Fixed Quarter Date'
$ws.Range('D204').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D226').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D281').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D291').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D297').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D299').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D317').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D322').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D378').Value = 'Downgrade Since Last Fixed Quarter Date'
$ws.Range('D406').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D407').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D414').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D425').Value = 'Downgrade Since Last Fixed Quarter Date'
$ws.Range('D429').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D465').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D475').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D486').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D497').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D506').Value = 'Downgrade Since Last Fixed Quarter Date'
$ws.Range('D545').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D555').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D565').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D568').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D591').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D598').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D599').Value = 'Downgrade Since Last Fixed Quarter Date'
$ws.Range('D622').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D636').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D643').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D652').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D688').Value = 'Downgrade Since Last Fixed Quarter Date'
$ws.Range('D711').Value = 'Downgrade Since Last Fixed Quarter Date'
$ws.Range('D717').Value = 'Downgrade Since Last Fixed Quarter Date'
$ws.Range('D742').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D743').Value = 'Downgrade Since Last Fixed Quarter Date'
$ws.Range('D749').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D764').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D765').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D801').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D886').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D965').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D996').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D1011').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D1016').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D1031').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D1051').Value = 'Downgrade Since Last Fixed Quarter Date'
$ws.Range('D1060').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D1068').Value = 'Downgrade Since Last Fixed Quarter Date'
$ws.Range('D1097').Value = 'Upgrade Since Last Fixed Quarter Date'
$ws.Range('D1116').Value = 'Same As Last Fixed Quarter Date'
$ws.Range('D1118').Value = 'Same As Last Fixed Quarter Date'
